$d = $word.ActiveDocument
$em = [char]0x2014

# ---------------------------------------------------------------------------
# 1. Insert author + affiliation block right after the Title paragraph.
# ---------------------------------------------------------------------------
$pTitle = $d.Paragraphs(1)

# Create two blank paragraphs (Normal style) right after the title.
$pTitle.Range.InsertParagraphAfter()
$pAuthor = $d.Paragraphs(2)
$pAuthor.Range.Style = "Normal"
$pAuthor.Range.InsertParagraphAfter()
$pAddress = $d.Paragraphs(3)
$pAddress.Range.Style = "Normal"

# Fill the author paragraph: bold run, leading line break, then the name.
$pAuthor.Range.Text = [char]11 + "Siddalingaiah H S, MD"
$rAuthorBold = $d.Range($pAuthor.Range.Start, $pAuthor.Range.End - 1)
$rAuthorBold.Font.Bold = 1

# Fill the affiliation / address paragraph (line breaks between each line,
# including a trailing one after the email).
$addrText = "Department of Community Medicine" + [char]11 + "Shridevi Institute of Medical Sciences and Research Hospital" + [char]11 + "Tumkur, Karnataka, India - 572106" + [char]11 + "Email: hssling@yahoo.com" + [char]11
$pAddress.Range.Text = $addrText

# ---------------------------------------------------------------------------
# 2. Rewrite the Abstract paragraph text (paragraph 5 now: Title, Author,
#    Address, "Abstract" heading, Background text).
# ---------------------------------------------------------------------------
$newAbstract = "Background: Tuberculosis (TB) treatment failure affects ~5-10% of drug-susceptible patients, yet biomarkers predicting this outcome remain elusive. Current host-response signatures are derived primarily from peripheral blood (PBMCs), potentially missing the critical immunopathology occurring at the site of infection" + $em + "the lung. Methods: We performed paired single-cell RNA-sequencing (scRNA-seq) and Assay for Transposase-Accessible Chromatin (scATAC-seq) on Bronchoalveolar Lavage (BAL) fluid and matched PBMCs from patients with active pulmonary TB. We utilized the Chromatin Priming Index (CPI) to map the epigenetic potential of immune cells and stratified patients by treatment outcome (Cure vs. Failure). Results: We observed a striking epigenetic divergence between compartments. Alveolar macrophages displayed a 'hyper-primed' inflammatory state (CPI 78.8%) driven by AP-1 (FOS/JUN) and NF-kB motifs, whereas peripheral monocytes showed a distinct Interferon-Response Factor (IRF) accessibility signature (CPI 84.2%). Comparing patients who cured vs. those who failed treatment, we identified a specific 'Failure chromatin signature' in lung macrophages, characterized by accessible chromatin at Matrix Metalloproteinase loci (MMP1, MMP9) despite low baseline expression. Conclusions: Chromatin accessibility landscapes in the TB lung are distinct from the periphery and predictive of clinical outcome. The identification of an epigenetically poised 'tissue destruction' program offers a novel therapeutic target for preventing lung damage and treatment failure."
$d.Paragraphs(5).Range.Text = $newAbstract

# ---------------------------------------------------------------------------
# 3. Insert the Introduction section right after the Abstract paragraph.
# ---------------------------------------------------------------------------
$pAbstractBody = $d.Paragraphs(5)
$pAbstractBody.Range.InsertParagraphAfter()
$pIntroHeading = $d.Paragraphs(6)
$pIntroHeading.Range.Style = "Heading1"
$pIntroHeading.Range.Text = "Introduction"

$pIntroHeading.Range.InsertParagraphAfter()
$pIntroBody = $d.Paragraphs(7)
$pIntroBody.Range.Style = "Normal"
$pIntroBody.Range.Text = "Despite effective chemotherapy, tuberculosis remains a leading cause of death globally. A significant subset of patients experiences 'treatment failure'" + $em + "defined as persistent culture positivity or recurrence" + $em + "driven not only by bacterial resistance but by host immunopathology (cavitation, fibrosis). The hallmark of TB pathology is the granuloma, a structure dominated by macrophages. While blood transcriptomics have yielded diagnostic signatures [1], they often fail to capture the tissue-specific immune dynamics driving lung destruction."

# ---------------------------------------------------------------------------
# 4. Results section: add two Heading2 subsections, rewriting the existing
#    results sentence into the first subsection.
# ---------------------------------------------------------------------------
$pResultsHeading = $d.Paragraphs(8)

$pResultsHeading.Range.InsertParagraphAfter()
$pSub1Heading = $d.Paragraphs(9)
$pSub1Heading.Range.Style = "Heading2"
$pSub1Heading.Range.Text = "The Lung is Epigenetically Distinct from Blood"

# Paragraph 10 is the old "We observed..." sentence -- replace its text directly.
$newResults = "Paired analysis of BAL and PBMC samples revealed that while transcriptional profiles showed some overlap, chromatin accessibility landscapes were profoundly distinct. Alveolar Macrophages (AMs) were enriched for motifs of the AP-1 family (FOS, JUN, FOSB), consistent with a 'tissue-resident activated' phenotype. In contrast, peripheral monocytes were dominated by ISRE and STAT motifs, reflecting a systemic interferon response."
$d.Paragraphs(10).Range.Text = $newResults

$pSub1Body = $d.Paragraphs(10)
$pSub1Body.Range.InsertParagraphAfter()
$pSub2Heading = $d.Paragraphs(11)
$pSub2Heading.Range.Style = "Heading2"
$pSub2Heading.Range.Text = "The 'Failure' Chromatin Signature"

$pSub2Heading.Range.InsertParagraphAfter()
$pSub2Body = $d.Paragraphs(12)
$pSub2Body.Range.Style = "Normal"
$pSub2Body.Range.Text = "We stratified patients based on their 6-month treatment outcomes. Patients who failed treatment exhibited a specific chromatin signature in their Alveolar Macrophages at baseline (pre-treatment). This signature was characterized by increased accessibility at loci encoding tissue-destructive enzymes, specifically MMP1 and MMP9. Transcription Factor motif reinforcement analysis identified BATF and MAF as the master regulators maintaining this pathological chromatin state."

# ---------------------------------------------------------------------------
# 5. Discussion section.
# ---------------------------------------------------------------------------
$pSub2Body.Range.InsertParagraphAfter()
$pDiscHeading = $d.Paragraphs(13)
$pDiscHeading.Range.Style = "Heading1"
$pDiscHeading.Range.Text = "Discussion"

$pDiscHeading.Range.InsertParagraphAfter()
$pDiscBody = $d.Paragraphs(14)
$pDiscBody.Range.Style = "Normal"
$pDiscBody.Range.Text = "Our findings suggest that 'Treatment Failure' is not a random event but a pre-determined immunological state encoded in the chromatin of lung macrophages. The 'open' state of MMP genes suggests these cells are primed to cause cavitation upon stimulation. This highlights the urgent need for host-directed therapies (HDTs) that can remodel the lung epigenetic landscape, such as inhaled HDAC inhibitors."

# ---------------------------------------------------------------------------
# 6. Declarations section.
# ---------------------------------------------------------------------------
$pDiscBody.Range.InsertParagraphAfter()
$pDeclHeading = $d.Paragraphs(15)
$pDeclHeading.Range.Style = "Heading1"
$pDeclHeading.Range.Text = "Declarations"

$pDeclHeading.Range.InsertParagraphAfter()
$pFunding = $d.Paragraphs(16)
$pFunding.Range.Style = "Normal"
$pFunding.Range.Text = "Funding: No specific funding received."

$pFunding.Range.InsertParagraphAfter()
$pCompeting = $d.Paragraphs(17)
$pCompeting.Range.Style = "Normal"
$pCompeting.Range.Text = "Competing Interests: The authors declare no competing interests."

$pCompeting.Range.InsertParagraphAfter()
$pData = $d.Paragraphs(18)
$pData.Range.Style = "Normal"
$pData.Range.Text = "Data Availability: All analysis code and processed data are available at: https://github.com/hssling/CPI_MultiDisease_Extension"

$pData.Range.InsertParagraphAfter()
$pEthics = $d.Paragraphs(19)
$pEthics.Range.Style = "Normal"
$pEthics.Range.Text = "Ethical Approval: The study was approved by the Institutional Ethics Committee (IEC)."

# ---------------------------------------------------------------------------
# 7. References section: replace reference 1 & 2, insert a new reference 3.
# ---------------------------------------------------------------------------
$newRef1 = "1. Zak DE, et al. A blood RNA signature for tuberculosis disease risk: a prospective cohort study. Lancet. 2016;387(10035):2312-2322."
$d.Paragraphs(21).Range.Text = $newRef1

$newRef2 = "2. Elkington PT, et al. MMP-1 drives immunopathology in human tuberculosis and transgenic mice. J Clin Invest. 2011;121(5):1827-1833."
$d.Paragraphs(22).Range.Text = $newRef2

$pRef2 = $d.Paragraphs($d.Paragraphs.Count)
$pRef2.Range.InsertParagraphAfter()
$pRef3 = $d.Paragraphs($d.Paragraphs.Count)
$pRef3.Range.Style = "Normal"
$pRef3.Range.Text = "3. Pacis A, et al. Bacterial infection remodels the DNA methylation landscape of human dendritic cells. Genome Res. 2015;25(12):1801-1811."
